$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 102
$ws.Cells.Item(2, 4).Value = 94.5
$ws.Cells.Item(3, 3).Value = 31
$ws.Cells.Item(3, 4).Value = 40.5
$ws.Cells.Item(4, 3).Value = 192
$ws.Cells.Item(4, 4).Value = 190.5
$ws.Cells.Item(5, 3).Value = 362
$ws.Cells.Item(5, 4).Value = 354
$ws.Cells.Item(6, 3).Value = 132
$ws.Cells.Item(6, 4).Value = 116.5
$ws.Cells.Item(7, 3).Value = 189
$ws.Cells.Item(7, 4).Value = 188.5
$ws.Cells.Item(8, 3).Value = 48
$ws.Cells.Item(8, 4).Value = 45.5
$ws.Cells.Item(9, 3).Value = 141
$ws.Cells.Item(9, 4).Value = 135
$ws.Cells.Item(11, 3).Value = 38
$ws.Cells.Item(11, 4).Value = 44.5
$ws.Cells.Item(12, 3).Value = 133
$ws.Cells.Item(12, 4).Value = 123.5
$ws.Cells.Item(13, 3).Value = 32
$ws.Cells.Item(13, 4).Value = 32.5
$ws.Cells.Item(14, 3).Value = 219
$ws.Cells.Item(14, 4).Value = 212
$ws.Cells.Item(15, 3).Value = 210
$ws.Cells.Item(15, 4).Value = 226.5
$ws.Cells.Item(16, 3).Value = 14
$ws.Cells.Item(16, 4).Value = 14
$ws.Cells.Item(17, 3).Value = 29
$ws.Cells.Item(17, 4).Value = 27
$ws.Cells.Item(18, 3).Value = 72
$ws.Cells.Item(18, 4).Value = 76.5
$ws.Cells.Item(19, 3).Value = 40
$ws.Cells.Item(19, 4).Value = 33.5
$ws.Cells.Item(20, 3).Value = 11
$ws.Cells.Item(20, 4).Value = 8
$ws.Cells.Item(21, 3).Value = 70
$ws.Cells.Item(21, 4).Value = 56
$ws.Cells.Item(22, 3).Value = 87
$ws.Cells.Item(22, 4).Value = 89.5
$ws.Cells.Item(23, 3).Value = 51
$ws.Cells.Item(23, 4).Value = 52.5
$ws.Cells.Item(24, 3).Value = 129
$ws.Cells.Item(24, 4).Value = 126.5
$ws.Cells.Item(25, 3).Value = 12
$ws.Cells.Item(25, 4).Value = 9
$ws.Cells.Item(26, 3).Value = 114
$ws.Cells.Item(26, 4).Value = 113.5
$ws.Cells.Item(27, 3).Value = 42
$ws.Cells.Item(27, 4).Value = 34
$ws.Cells.Item(28, 3).Value = 71
$ws.Cells.Item(28, 4).Value = 58
$ws.Cells.Item(29, 3).Value = 56
$ws.Cells.Item(29, 4).Value = 57.5
$ws.Cells.Item(30, 3).Value = 78
$ws.Cells.Item(30, 4).Value = 62
$ws.Cells.Item(31, 3).Value = 10
$ws.Cells.Item(31, 4).Value = 8.5
$ws.Cells.Item(32, 3).Value = 344
$ws.Cells.Item(32, 4).Value = 363
$ws.Cells.Item(33, 3).Value = 172
$ws.Cells.Item(33, 4).Value = 166
$ws.Cells.Item(34, 3).Value = 26
$ws.Cells.Item(34, 4).Value = 24.5
$ws.Cells.Item(35, 3).Value = 19
$ws.Cells.Item(35, 4).Value = 14.5
$ws.Cells.Item(36, 3).Value = 101
$ws.Cells.Item(36, 4).Value = 92.5
$ws.Cells.Item(37, 3).Value = 43
$ws.Cells.Item(37, 4).Value = 36
$ws.Cells.Item(38, 3).Value = 35
$ws.Cells.Item(38, 4).Value = 25
$ws.Cells.Item(39, 3).Value = 95.02702702702703
